{"js": "// Insert a new paragraph \"what is your problem brother\" at the very start\n// of the document body (before whatever paragraph is currently first),\n// reproducing Word's grammar-checker markup: the word \"brother\" (a comma\n// splice / direct-address flag) is wrapped in a gramStart/gramEnd\n// w:proofErr pair, which also splits the run in two.\nconst body = context.document.body;\n\n// Find the body's current first paragraph so the new one can be inserted\n// immediately \"Before\" it (falls back to inserting at the body's End if the\n// body is somehow empty of paragraphs).\nconst firstParagraph = body.paragraphs.getFirstOrNullObject();\nawait context.sync();\n\nconst newParagraph = firstParagraph.isNullObject\n  ? body.insertParagraph(\"\", \"End\")\n  : firstParagraph.insertParagraph(\"\", \"Before\");\nawait context.sync();\n\n// Flat-OPC wrapped OOXML for the paragraph's final contents: two runs\n// (\"what is your problem \" / \"brother\") separated by the proofErr pair,\n// exactly matching what Word's proofing pass emits for this sentence.\nconst paragraphOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">what is your problem </w:t></w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r><w:t>brother</w:t></w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nnewParagraph.insertOoxml(paragraphOoxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Insert a new paragraph \"what is your problem brother\" at the very start\n# of the document body (before whatever paragraph is currently first),\n# reproducing Word's grammar-checker markup: the word \"brother\" (flagged as\n# a direct-address comma splice) is wrapped in a gramStart/gramEnd\n# w:proofErr pair, which also splits the run in two.\n$d = $word.ActiveDocument\n\n# Insert a fresh empty paragraph immediately before the document's current\n# first paragraph, leaving that original paragraph completely untouched.\n$firstParagraph = $d.Paragraphs.First\n$firstParagraph.Range.InsertParagraphBefore() | Out-Null\n\n# The freshly inserted paragraph is now the document's first paragraph;\n# fill it in via raw WordOpenXML so the run split + proofErr markers come\n# out exactly as Word's proofing pass would write them.\n$newParagraph = $d.Paragraphs.First\n$paragraphXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t xml:space=\"preserve\">what is your problem </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>brother</w:t></w:r><w:proofErr w:type=\"gramEnd\"/></w:p>'\n$newParagraph.Range.InsertXML($paragraphXml) | Out-Null\n"}
